# Update the cryptos worksheet with the latest scraped price/volume figures.
# Two coin pairs (Dai/ShibaInu and WrappedBTC/Toncoin) swapped rank order,
# so their Coin name + Link values are swapped along with fresh price/volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row swaps: Coin name (B) and Link (C) exchange places between the two rows ---

# Rows 16 <-> 17 : Dai / ShibaInu
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"

# Rows 24 <-> 25 : WrappedBTC / Toncoin
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"

# --- Price (D) and Volume(1h) (E) updates for every data row (2-51, row 4 unchanged) ---
# D/E are stored as plain text in the sheet (e.g. "21.761.97", "  +5.98%  "), so force
# the cells to Text format first to stop Excel from re-interpreting/reformatting the
# numeric-looking strings (which would strip trailing zeros / change precision).

$updates = @{
    2  = @{ D = "21.761.97";  E = "  +5.98%  " }
    3  = @{ D = "1.572.50";   E = "  +6.80%  " }
    5  = @{ D = "0.9844";     E = "  +2.85%  " }
    6  = @{ D = "285.80";     E = "  +3.46%  " }
    7  = @{ D = "0.3686";     E = "  +1.05%  " }
    8  = @{ D = "0.3277";     E = "  +7.04%  " }
    9  = @{ D = "41.78";      E = "  +4.98%  " }
    10 = @{ D = "1.134";      E = "  +7.47%  " }
    11 = @{ D = "0.07041";    E = "  +6.34%  " }
    12 = @{ D = "0.9975";     E = "  -0.40%  " }
    13 = @{ D = "19.92";      E = "  +10.00%  " }
    14 = @{ D = "5.836";      E = "  +6.75%  " }
    15 = @{ D = "6.514";      E = "  +5.40%  " }
    16 = @{ D = "0.00001070"; E = "  +3.99%  " }
    17 = @{ D = "0.9844";     E = "  +2.15%  " }
    18 = @{ D = "1.565.91";   E = "  +6.25%  " }
    19 = @{ D = "0.06215";    E = "  +5.43%  " }
    20 = @{ D = "74.55";      E = "  +7.80%  " }
    21 = @{ D = "16.19";      E = "  +11.52%  " }
    22 = @{ D = "5.866";      E = "  +7.31%  " }
    23 = @{ D = "11.60";      E = "  +4.97%  " }
    24 = @{ D = "2.368";      E = "  +5.38%  " }
    25 = @{ D = "21.740.12";  E = "  +5.60%  " }
    26 = @{ D = "2.370";      E = "  +11.35%  " }
    27 = @{ D = "149.31";     E = "  +6.38%  " }
    28 = @{ D = "18.25";      E = "  +6.09%  " }
    29 = @{ D = "1.740.44";   E = "  +6.71%  " }
    30 = @{ D = "120.48";     E = "  +5.70%  " }
    31 = @{ D = "4.093";      E = "  +2.87%  " }
    32 = @{ D = "0.9084";     E = "  +11.68%  " }
    33 = @{ D = "5.435";      E = "  +9.50%  " }
    34 = @{ D = "0.08197";    E = "  +3.21%  " }
    35 = @{ D = "1.597";      E = "  +4.06%  " }
    36 = @{ D = "5.117";      E = "  +8.66%  " }
    37 = @{ D = "1.238";      E = "  +1.21%  " }
    38 = @{ D = "11.51";      E = "  +10.41%  " }
    39 = @{ D = "0.06043";    E = "  +4.26%  " }
    40 = @{ D = "0.02170";    E = "  +6.50%  " }
    41 = @{ D = "8.126";      E = "  +7.00%  " }
    42 = @{ D = "0.1994";     E = "  +6.17%  " }
    43 = @{ D = "0.9841";     E = "  +2.62%  " }
    44 = @{ D = "0.5742";     E = "  +8.97%  " }
    45 = @{ D = "12.94";      E = "  +7.60%  " }
    46 = @{ D = "3.619";      E = "  +3.15%  " }
    47 = @{ D = "0.5622";     E = "  +8.63%  " }
    48 = @{ D = "124.60";     E = "  +5.84%  " }
    49 = @{ D = "1.912";      E = "  +6.87%  " }
    50 = @{ D = "0.06737";    E = "  +4.28%  " }
    51 = @{ D = "71.76";      E = "  +7.14%  " }
}

foreach ($rowNum in $updates.Keys) {
    $vals = $updates[$rowNum]

    $dCell = $ws.Range("D$rowNum")
    $dCell.NumberFormat = "@"
    $dCell.Value = $vals.D

    $eCell = $ws.Range("E$rowNum")
    $eCell.NumberFormat = "@"
    $eCell.Value = $vals.E
}
